$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update masked ID cells (C2:C4) from 7-char mask to 6-char mask
$ws.Range("C2").Value = "######"
$ws.Range("C3").Value = "######"
$ws.Range("C4").Value = "######"

# Update similarity test cells (D2:D4) from 4-char to 6-char mask
$ws.Range("D2").Value = "&&&&&&"
$ws.Range("D3").Value = "&&&&&&"
$ws.Range("D4").Value = "&&&&&&"

# Update selection to H7
$ws.Range("H7").Select()
